# Automatische test-sync: 2025-08-04 19:59:50
# Adds a new log row (row 3) to the "Logs" sheet, extends the conditional
# formatting ranges that applied only to row 2 so they also cover row 3,
# and bumps the "Aantal" count on the "Dashboard" sheet from 1 to 2.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #1: Hoi, hebben jullie al iets gehoord?"
$logs.Range("D3").Value = "Opvolging / Status"
$logs.Range("E3").Value = "Dank voor je bericht. We hebben je eerdere e-mail ontvangen en doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Range("F3").Value = "2025-08-04 19:59:12"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to include row 3 --
# Each of these columns had a conditionalFormatting block scoped to the
# single cell in row 2; stretch every rule in each block down to row 3.

$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldCell = "$col" + "2"
    $newRange = "$col" + "2:$col" + "3"
    $fcs = $logs.Range($oldCell).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Dashboard sheet: bump the count -------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 2
